$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 41275.6
$ws.Range("I33").Value = 45904.09
$ws.Range("K33").Value = 45904.09
$ws.Range("M33").Value = -45675.09
$ws.Range("H64").Value = 128012.5
$ws.Range("I64").Value = 335500
$ws.Range("J64").Value = 3520
$ws.Range("K64").Value = 335500
$ws.Range("L64").Value = 3520
$ws.Range("M64").Value = -335252
$ws.Range("N64").Value = -4016
$ws.Range("H67").Value = 128012.5
$ws.Range("I67").Value = 335500
$ws.Range("J67").Value = 3520
$ws.Range("K67").Value = 335500
$ws.Range("L67").Value = 3520
$ws.Range("M67").Value = -334642
$ws.Range("N67").Value = -5236
$ws.Range("H76").Value = 4487.2354
$ws.Range("I76").Value = 3948.125
$ws.Range("J76").Value = 4966.4443
$ws.Range("K76").Value = 3948.125
$ws.Range("L76").Value = 4966.4443
$ws.Range("M76").Value = -3633.125
$ws.Range("N76").Value = -5596.4443
$ws.Range("H79").Value = 4487.2354
$ws.Range("I79").Value = 3948.125
$ws.Range("J79").Value = 4966.4443
$ws.Range("K79").Value = 3948.125
$ws.Range("L79").Value = 4966.4443
$ws.Range("M79").Value = -2856.125
$ws.Range("N79").Value = -7150.4443
$ws.Range("H96").Value = 824.9
$ws.Range("I96").Value = 590
$ws.Range("J96").Value = 925.5714
$ws.Range("K96").Value = 1770
$ws.Range("L96").Value = 2776.7142
$ws.Range("M96").Value = -397
$ws.Range("N96").Value = -5522.7142
$ws.Range("H100").Value = 963.4
$ws.Range("I100").Value = 825
$ws.Range("J100").Value = 1032.6
$ws.Range("K100").Value = 825
$ws.Range("L100").Value = 1032.6
$ws.Range("M100").Value = -284
$ws.Range("N100").Value = -2114.6
$ws.Range("H112").Value = 1563407.5
$ws.Range("J112").Value = 1645660.5
$ws.Range("L112").Value = 4936981.5
$ws.Range("N112").Value = -4939197.5
$ws.Range("H135").Value = 1280.2903
$ws.Range("I135").Value = 720.5833
$ws.Range("J135").Value = 3199.2856
$ws.Range("K135").Value = 6485.2497
$ws.Range("L135").Value = 28793.5704
$ws.Range("M135").Value = -3950.2497
$ws.Range("N135").Value = -33863.5704
$ws.Range("H137").Value = 1808.4872
$ws.Range("I137").Value = 1220.3235
$ws.Range("K137").Value = 3660.9705
$ws.Range("M137").Value = -1110.9705

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 6214.2856
$ws.Range("J3").Value = 8000
$ws.Range("L3").Value = 8000
$ws.Range("N3").Value = -8230
$ws.Range("H32").Value = 24368.162
$ws.Range("I32").Value = 4003.3333
$ws.Range("J32").Value = 305402.8
$ws.Range("K32").Value = 4003.3333
$ws.Range("L32").Value = 305402.8
$ws.Range("M32").Value = -3716.3333
$ws.Range("N32").Value = -305976.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 64382.945
$ws.Range("I86").Value = 124545.11
$ws.Range("J86").Value = 4220.778
$ws.Range("K86").Value = 124545.11
$ws.Range("L86").Value = 4220.778
$ws.Range("M86").Value = -123422.11
$ws.Range("N86").Value = -6466.778
$ws.Range("H89").Value = 64382.945
$ws.Range("I89").Value = 124545.11
$ws.Range("J89").Value = 4220.778
$ws.Range("K89").Value = 622725.55
$ws.Range("L89").Value = 21103.89
$ws.Range("M89").Value = -617109.55
$ws.Range("N89").Value = -32335.89
$ws.Range("H94").Value = 519.28
$ws.Range("I94").Value = 384.5
$ws.Range("J94").Value = 758.8889
$ws.Range("K94").Value = 384.5
$ws.Range("L94").Value = 758.8889
$ws.Range("M94").Value = 66.5
$ws.Range("N94").Value = -1660.8889
$ws.Range("H99").Value = 1754.1875
$ws.Range("I99").Value = 1858
$ws.Range("J99").Value = 1707
$ws.Range("K99").Value = 1858
$ws.Range("L99").Value = 1707
$ws.Range("M99").Value = -360
$ws.Range("N99").Value = -4703
$ws.Range("H105").Value = 169974.25
$ws.Range("I105").Value = 114777.78
$ws.Range("J105").Value = 335563.66
$ws.Range("K105").Value = 114777.78
$ws.Range("L105").Value = 335563.66
$ws.Range("M105").Value = -113030.78
$ws.Range("N105").Value = -339057.66
$ws.Range("H134").Value = 2811.9556
$ws.Range("I134").Value = 3031.8064
$ws.Range("J134").Value = 2325.1428
$ws.Range("K134").Value = 9095.4192
$ws.Range("L134").Value = 6975.428400000001
$ws.Range("M134").Value = -6560.4192
$ws.Range("N134").Value = -12045.4284

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 11500
$ws.Range("I4").Value = 8000
$ws.Range("J4").Value = 13833.333
$ws.Range("K4").Value = 8000
$ws.Range("L4").Value = 13833.333
$ws.Range("M4").Value = -7888
$ws.Range("N4").Value = -14057.333
$ws.Range("H31").Value = 34061.266
$ws.Range("I31").Value = 1121.9445
$ws.Range("J31").Value = 53187.324
$ws.Range("K31").Value = 1121.9445
$ws.Range("L31").Value = 53187.324
$ws.Range("M31").Value = -826.9445000000001
$ws.Range("N31").Value = -53777.324
$ws.Range("H34").Value = 34061.266
$ws.Range("I34").Value = 1121.9445
$ws.Range("J34").Value = 53187.324
$ws.Range("K34").Value = 1121.9445
$ws.Range("L34").Value = 53187.324
$ws.Range("M34").Value = -919.9445000000001
$ws.Range("N34").Value = -53591.324
$ws.Range("H58").Value = 4942.244
$ws.Range("I58").Value = 1181.8438
$ws.Range("J58").Value = 18312.555
$ws.Range("K58").Value = 1181.8438
$ws.Range("L58").Value = 18312.555
$ws.Range("M58").Value = -978.8438000000001
$ws.Range("N58").Value = -18718.555
$ws.Range("H107").Value = 8436.385
$ws.Range("I107").Value = 15229.857
$ws.Range("J107").Value = 510.66666
$ws.Range("K107").Value = 15229.857
$ws.Range("L107").Value = 510.66666
$ws.Range("M107").Value = -13309.857
$ws.Range("N107").Value = -4350.66666
$ws.Range("H134").Value = 1353.8667
$ws.Range("I134").Value = 1313
$ws.Range("J134").Value = 1400.5714
$ws.Range("K134").Value = 3939
$ws.Range("L134").Value = 4201.7142
$ws.Range("M134").Value = -1404
$ws.Range("N134").Value = -9271.7142
$ws.Range("H136").Value = 4942.244
$ws.Range("I136").Value = 1181.8438
$ws.Range("J136").Value = 18312.555
$ws.Range("K136").Value = 3545.5314
$ws.Range("L136").Value = 54937.665
$ws.Range("M136").Value = -995.5314000000003
$ws.Range("N136").Value = -60037.665

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 20234.375
$ws.Range("J32").Value = 22839.285
$ws.Range("L32").Value = 68517.855
$ws.Range("N32").Value = -69083.855
$ws.Range("H140").Value = 1734.5483
$ws.Range("I140").Value = 1505.2778
$ws.Range("J140").Value = 2052
$ws.Range("K140").Value = 4515.8334
$ws.Range("L140").Value = 6156
$ws.Range("M140").Value = 664.1665999999996
$ws.Range("N140").Value = -16516
$ws.Range("H141").Value = 4500
$ws.Range("I141").Value = 3000
$ws.Range("K141").Value = 9000
$ws.Range("M141").Value = -3820

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H105").Value = 46325
$ws.Range("J105").Value = 46325
$ws.Range("L105").Value = 46325
$ws.Range("N105").Value = -53313
$ws.Range("H122").Value = 3902
$ws.Range("I122").Value = 6000
$ws.Range("J122").Value = 1804
$ws.Range("K122").Value = 18000
$ws.Range("L122").Value = 5412
$ws.Range("M122").Value = -15550
$ws.Range("N122").Value = -10312

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 723665
$ws.Range("I46").Value = 419.6
$ws.Range("J46").Value = 1125468
$ws.Range("K46").Value = 419.6
$ws.Range("L46").Value = 1125468
$ws.Range("M46").Value = -231.6
$ws.Range("N46").Value = -1125844
$ws.Range("H55").Value = 382.14285
$ws.Range("I55").Value = 231.23077
$ws.Range("J55").Value = 512.93335
$ws.Range("K55").Value = 231.23077
$ws.Range("L55").Value = 512.93335
$ws.Range("M55").Value = -58.23077000000001
$ws.Range("N55").Value = -858.93335
$ws.Range("H68").Value = 2847
$ws.Range("I68").Value = 1459.8
$ws.Range("J68").Value = 3380.5386
$ws.Range("K68").Value = 1459.8
$ws.Range("L68").Value = 3380.5386
$ws.Range("M68").Value = -710.8
$ws.Range("N68").Value = -4878.5386
$ws.Range("H71").Value = 2847
$ws.Range("I71").Value = 1459.8
$ws.Range("J71").Value = 3380.5386
$ws.Range("K71").Value = 7299
$ws.Range("L71").Value = 16902.693
$ws.Range("M71").Value = -3555
$ws.Range("N71").Value = -24390.693
$ws.Range("H132").Value = 2625.9565
$ws.Range("J132").Value = 2244.125
$ws.Range("L132").Value = 6732.375
$ws.Range("N132").Value = -11792.375
$ws.Range("H136").Value = 2070.5
$ws.Range("I136").Value = 1889.2222
$ws.Range("K136").Value = 5667.6666
$ws.Range("M136").Value = -3117.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 251122.38
$ws.Range("J81").Value = 251444.75
$ws.Range("L81").Value = 502889.5
$ws.Range("N81").Value = -505011.5
$ws.Range("H84").Value = 251122.38
$ws.Range("J84").Value = 251444.75
$ws.Range("L84").Value = 2514447.5
$ws.Range("N84").Value = -2525055.5
$ws.Range("H132").Value = 1982.5
$ws.Range("I132").Value = 1895.4348
$ws.Range("J132").Value = 2316.25
$ws.Range("K132").Value = 5686.3044
$ws.Range("L132").Value = 6948.75
$ws.Range("M132").Value = -3156.3044
$ws.Range("N132").Value = -12008.75
$ws.Range("H136").Value = 1400.8334
$ws.Range("I136").Value = 921
$ws.Range("J136").Value = 2000.625
$ws.Range("K136").Value = 2763
$ws.Range("L136").Value = 6001.875
$ws.Range("M136").Value = -213
$ws.Range("N136").Value = -11101.875
